# Scenario 1 / 2 / 3 / 7 modifications - "01 - Registracija vlasnika.docx"
#
# 1) Table header cells "NAZIV" / "OPIS" get " SLUČAJA UPOTREBE" appended
#    (bold, same as the rest of the cell text).
# 2) "Uspješno registrovan novi vlasnik" - the two runs that made up this
#    sentence are unified into a single run with the same text.
# 3) The "USPJEŠAN ZAVRŠETAK" / "NEUSPJEŠAN ZAVRŠETAK" headings: the
#    "<name> " + "–" + " " run triplet becomes a single "<name> – " run.
# 4) "Unos potrebnih podataka" (appears twice) - the three runs are
#    unified into a single run with the same text.

$d = $word.ActiveDocument

# --- 1a) NAZIV -> NAZIV SLUČAJA UPOTREBE ------------------------------
$d.Content.Find.Execute(
    "NAZIV", $true, $false, $false, $false, $false, $true, 1, $false,
    "NAZIV SLUČAJA UPOTREBE", 2) | Out-Null

# --- 1b) OPIS -> OPIS SLUČAJA UPOTREBE --------------------------------
$d.Content.Find.Execute(
    "OPIS", $true, $false, $false, $false, $false, $true, 1, $false,
    "OPIS SLUČAJA UPOTREBE", 2) | Out-Null

# --- 2) "Uspješno " + "registrovan novi vlasnik" -> one run ----------
$d.Content.Find.Execute(
    "Uspješno registrovan novi vlasnik", $true, $false, $false, $false,
    $false, $true, 1, $false, "Uspješno registrovan novi vlasnik", 2) | Out-Null

# --- 3) Success / failure heading dash runs merged --------------------
$d.Content.Find.Execute(
    "USPJEŠAN ZAVRŠETAK " + [char]0x2013 + " ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "USPJEŠAN ZAVRŠETAK " + [char]0x2013 + " ", 2) | Out-Null

$d.Content.Find.Execute(
    "NEUSPJEŠAN ZAVRŠETAK " + [char]0x2013 + " ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "NEUSPJEŠAN ZAVRŠETAK " + [char]0x2013 + " ", 2) | Out-Null

# --- 4) "Unos " + "potrebnih" + " podataka" -> one run (both tables) --
$d.Content.Find.Execute(
    "Unos potrebnih podataka", $true, $false, $false, $false, $false,
    $true, 1, $false, "Unos potrebnih podataka", 2) | Out-Null
